$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (styles) for new rows 2098:2196 by copying format from 2067:2097
$ws.Range("A2067:B2097").Copy()
$ws.Range("A2098:B2196").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update revised values for rows 2067-2097 (Dec 2024 re-aggregation)
$ws.Range("B2067").Value = 47.694
$ws.Range("B2068").Value = 53.864
$ws.Range("B2069").Value = 50.821
$ws.Range("B2070").Value = 47.48
$ws.Range("B2071").Value = 43.846
$ws.Range("B2072").Value = 55.304
$ws.Range("B2073").Value = 46.862
$ws.Range("B2074").Value = 50.054
$ws.Range("B2075").Value = 40.466
$ws.Range("B2076").Value = 40.006
$ws.Range("B2077").Value = 42.207
$ws.Range("B2078").Value = 43.915
$ws.Range("B2079").Value = 40.317
$ws.Range("B2080").Value = 36.785
$ws.Range("B2081").Value = 40.831
$ws.Range("B2082").Value = 38.99
$ws.Range("B2083").Value = 40.31
$ws.Range("B2084").Value = 36.884
$ws.Range("B2085").Value = 37.57
$ws.Range("B2086").Value = 35.561
$ws.Range("B2087").Value = 35.078
$ws.Range("B2088").Value = 36.428
$ws.Range("B2089").Value = 33.32
$ws.Range("B2090").Value = 26.841
$ws.Range("B2091").Value = 31.39
$ws.Range("B2092").Value = 38.409
$ws.Range("B2093").Value = 34.167
$ws.Range("B2094").Value = 39.211
$ws.Range("B2095").Value = 36.293
$ws.Range("B2096").Value = 38.897
$ws.Range("B2097").Value = 29.829

# Append new rows 2098-2196 for Q1 2025 (Jan 1 - Apr 9, 2025)
$ws.Range("A2098").Value = 45658
$ws.Range("B2098").Value = 23.942
$ws.Range("A2099").Value = 45659
$ws.Range("B2099").Value = 31.733
$ws.Range("A2100").Value = 45660
$ws.Range("B2100").Value = 34.336
$ws.Range("A2101").Value = 45661
$ws.Range("B2101").Value = 36.812
$ws.Range("A2102").Value = 45662
$ws.Range("B2102").Value = 32.368
$ws.Range("A2103").Value = 45663
$ws.Range("B2103").Value = 32.813
$ws.Range("A2104").Value = 45664
$ws.Range("B2104").Value = 33.518
$ws.Range("A2105").Value = 45665
$ws.Range("B2105").Value = 39.728
$ws.Range("A2106").Value = 45666
$ws.Range("B2106").Value = 32.077
$ws.Range("A2107").Value = 45667
$ws.Range("B2107").Value = 33.413
$ws.Range("A2108").Value = 45668
$ws.Range("B2108").Value = 32.253
$ws.Range("A2109").Value = 45669
$ws.Range("B2109").Value = 30.404
$ws.Range("A2110").Value = 45670
$ws.Range("B2110").Value = 31.513
$ws.Range("A2111").Value = 45671
$ws.Range("B2111").Value = 34.922
$ws.Range("A2112").Value = 45672
$ws.Range("B2112").Value = 39.248
$ws.Range("A2113").Value = 45673
$ws.Range("B2113").Value = 32.865
$ws.Range("A2114").Value = 45674
$ws.Range("B2114").Value = 34.599
$ws.Range("A2115").Value = 45675
$ws.Range("B2115").Value = 33.644
$ws.Range("A2116").Value = 45676
$ws.Range("B2116").Value = 33.843
$ws.Range("A2117").Value = 45677
$ws.Range("B2117").Value = 34.416
$ws.Range("A2118").Value = 45678
$ws.Range("B2118").Value = 32.595
$ws.Range("A2119").Value = 45679
$ws.Range("B2119").Value = 37.442
$ws.Range("A2120").Value = 45680
$ws.Range("B2120").Value = 33.333
$ws.Range("A2121").Value = 45681
$ws.Range("B2121").Value = 31.886
$ws.Range("A2122").Value = 45682
$ws.Range("B2122").Value = 32.49
$ws.Range("A2123").Value = 45683
$ws.Range("B2123").Value = 35.934
$ws.Range("A2124").Value = 45684
$ws.Range("B2124").Value = 32.494
$ws.Range("A2125").Value = 45685
$ws.Range("B2125").Value = 34.849
$ws.Range("A2126").Value = 45686
$ws.Range("B2126").Value = 30.919
$ws.Range("A2127").Value = 45687
$ws.Range("B2127").Value = 32.437
$ws.Range("A2128").Value = 45688
$ws.Range("B2128").Value = 30.15
$ws.Range("A2129").Value = 45689
$ws.Range("B2129").Value = 36.502
$ws.Range("A2130").Value = 45690
$ws.Range("B2130").Value = 29.632
$ws.Range("A2131").Value = 45691
$ws.Range("B2131").Value = 33.588
$ws.Range("A2132").Value = 45692
$ws.Range("B2132").Value = 36.163
$ws.Range("A2133").Value = 45693
$ws.Range("B2133").Value = 35.768
$ws.Range("A2134").Value = 45694
$ws.Range("B2134").Value = 35.43
$ws.Range("A2135").Value = 45695
$ws.Range("B2135").Value = 30.286
$ws.Range("A2136").Value = 45696
$ws.Range("B2136").Value = 41.019
$ws.Range("A2137").Value = 45697
$ws.Range("B2137").Value = 35.456
$ws.Range("A2138").Value = 45698
$ws.Range("B2138").Value = 37.344
$ws.Range("A2139").Value = 45699
$ws.Range("B2139").Value = 34.342
$ws.Range("A2140").Value = 45700
$ws.Range("B2140").Value = 29.718
$ws.Range("A2141").Value = 45701
$ws.Range("B2141").Value = 30.942
$ws.Range("A2142").Value = 45702
$ws.Range("B2142").Value = 28.804
$ws.Range("A2143").Value = 45703
$ws.Range("B2143").Value = 30.295
$ws.Range("A2144").Value = 45704
$ws.Range("B2144").Value = 31.785
$ws.Range("A2145").Value = 45705
$ws.Range("B2145").Value = 35.765
$ws.Range("A2146").Value = 45706
$ws.Range("B2146").Value = 36.331
$ws.Range("A2147").Value = 45707
$ws.Range("B2147").Value = 37.331
$ws.Range("A2148").Value = 45708
$ws.Range("B2148").Value = 51.691
$ws.Range("A2149").Value = 45709
$ws.Range("B2149").Value = 37.85
$ws.Range("A2150").Value = 45710
$ws.Range("B2150").Value = 37.989
$ws.Range("A2151").Value = 45711
$ws.Range("B2151").Value = 32.676
$ws.Range("A2152").Value = 45712
$ws.Range("B2152").Value = 35.736
$ws.Range("A2153").Value = 45713
$ws.Range("B2153").Value = 33.638
$ws.Range("A2154").Value = 45714
$ws.Range("B2154").Value = 33.861
$ws.Range("A2155").Value = 45715
$ws.Range("B2155").Value = 40.804
$ws.Range("A2156").Value = 45716
$ws.Range("B2156").Value = 26.375
$ws.Range("A2157").Value = 45717
$ws.Range("B2157").Value = 30.652
$ws.Range("A2158").Value = 45718
$ws.Range("B2158").Value = 28.622
$ws.Range("A2159").Value = 45719
$ws.Range("B2159").Value = 30.337
$ws.Range("A2160").Value = 45720
$ws.Range("B2160").Value = 33.669
$ws.Range("A2161").Value = 45721
$ws.Range("B2161").Value = 32.527
$ws.Range("A2162").Value = 45722
$ws.Range("B2162").Value = 30.17
$ws.Range("A2163").Value = 45723
$ws.Range("B2163").Value = 28.047
$ws.Range("A2164").Value = 45724
$ws.Range("B2164").Value = 33.301
$ws.Range("A2165").Value = 45725
$ws.Range("B2165").Value = 34.372
$ws.Range("A2166").Value = 45726
$ws.Range("B2166").Value = 36.206
$ws.Range("A2167").Value = 45727
$ws.Range("B2167").Value = 30.621
$ws.Range("A2168").Value = 45728
$ws.Range("B2168").Value = 30.753
$ws.Range("A2169").Value = 45729
$ws.Range("B2169").Value = 36.825
$ws.Range("A2170").Value = 45730
$ws.Range("B2170").Value = 31.924
$ws.Range("A2171").Value = 45731
$ws.Range("B2171").Value = 31.504
$ws.Range("A2172").Value = 45732
$ws.Range("B2172").Value = 35.167
$ws.Range("A2173").Value = 45733
$ws.Range("B2173").Value = 35.583
$ws.Range("A2174").Value = 45734
$ws.Range("B2174").Value = 41.897
$ws.Range("A2175").Value = 45735
$ws.Range("B2175").Value = 39.894
$ws.Range("A2176").Value = 45736
$ws.Range("B2176").Value = 38.547
$ws.Range("A2177").Value = 45737
$ws.Range("B2177").Value = 31.507
$ws.Range("A2178").Value = 45738
$ws.Range("B2178").Value = 32.219
$ws.Range("A2179").Value = 45739
$ws.Range("B2179").Value = 37.398
$ws.Range("A2180").Value = 45740
$ws.Range("B2180").Value = 38.186
$ws.Range("A2181").Value = 45741
$ws.Range("B2181").Value = 32.443
$ws.Range("A2182").Value = 45742
$ws.Range("B2182").Value = 34.434
$ws.Range("A2183").Value = 45743
$ws.Range("B2183").Value = 34.309
$ws.Range("A2184").Value = 45744
$ws.Range("B2184").Value = 30.534
$ws.Range("A2185").Value = 45745
$ws.Range("B2185").Value = 31.523
$ws.Range("A2186").Value = 45746
$ws.Range("B2186").Value = 32.909
$ws.Range("A2187").Value = 45747
$ws.Range("B2187").Value = 37.604
$ws.Range("A2188").Value = 45748
$ws.Range("B2188").Value = 29.867
$ws.Range("A2189").Value = 45749
$ws.Range("B2189").Value = 33.766
$ws.Range("A2190").Value = 45750
$ws.Range("B2190").Value = 47.498
$ws.Range("A2191").Value = 45751
$ws.Range("B2191").Value = 28.094
$ws.Range("A2192").Value = 45752
$ws.Range("B2192").Value = 33.214
$ws.Range("A2193").Value = 45753
$ws.Range("B2193").Value = 29.373
$ws.Range("A2194").Value = 45754
$ws.Range("B2194").Value = 39.574
$ws.Range("A2195").Value = 45755
$ws.Range("B2195").Value = 41.936
$ws.Range("A2196").Value = 45756
$ws.Range("B2196").Value = 30.515
